# "C+- corrigido e backup taylor 2GHz sincronizado"
#
# 1) Fix the mislabeled "C+/C-" combination-count annotation in the
#    RN_taylor sheet (cell K23) from "21(12)" to "21(4)".
# 2) Re-sync the Taylor-series backup columns (L = Elementos Logicos,
#    M = Multiplicadores, N = Memoria): only the 2 GHz row (row 23,
#    B23 = 2000 MHz) still has valid fresh backup data, so the older
#    rows (3-22) are cleared out, and the 2 GHz row is updated with
#    the newly synced values.

$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("RN_taylor")

# --- 1) fix the "21(12)" -> "21(4)" label -------------------------------
$ws.Range("K23").Value2 = "21(4)"

# --- 2) clear the stale backup rows and sync the 2 GHz row -------------
$ws.Range("L3:N22").ClearContents()

$ws.Range("L23").Value2 = 23810
$ws.Range("M23").Value2 = 112
$ws.Range("N23").Value2 = 259280

# Leave the cursor where the author left it when saving.
$ws.Activate()
$ws.Range("H2").Select()
